# Adds an explanatory note ("comment") after the source-code link at the
# end of the document: a blank line followed by a new paragraph reading
# "*Adicionar 5% no volume da fossa".
#
# Matches the commit: "adcionando comentarios ao docx"

$d = $word.ActiveDocument

# Locate the paragraph that holds the repository link ("...calc-bio").
# The two new paragraphs are inserted right after it, before the final
# (empty, italic, 10pt) paragraph that closes the document.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*calc-bio*") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    Write-Output "ERROR: anchor paragraph (containing 'calc-bio') not found"
} else {
    # 1) Split right after the anchor paragraph -> creates a new blank
    #    paragraph that inherits the anchor's paragraph formatting
    #    (justified, Times New Roman 12pt).
    $anchor.Range.InsertParagraphAfter()
    $blank = $anchor.Next()

    # 2) Split again after the (still empty) blank paragraph -> creates
    #    the paragraph that will hold the note text, with the same
    #    formatting.
    $blank.Range.InsertParagraphAfter()
    $note = $blank.Next()

    # 3) Fill in the note text.
    $note.Range.Text = "*Adicionar 5% no volume da fossa"

    Write-Output "Inserted note paragraph after '$($anchor.Range.Text.Trim())'"
}
